$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assign the same text to A1; this causes a new shared-string entry
# ("HK_R_acc_SD") to be inserted, giving a better (repeated/refreshed)
# description of the localT2 ls column header.
$ws.Range("A1").Value = "HK_R_acc_SD"

$newValues = @(
  78.057700121901661,
  77.651361235270215,
  78.017066233238523,
  79.926859000406338,
  78.179601787891102,
  80.414465664364073,
  75.822836245428689,
  77.895164567249083,
  73.872409589597723,
  76.757415684681021,
  78.13896789922795,
  80.333197887037784,
  71.271840715156443,
  72.247054043071927,
  72.003250711093045,
  88.90694839496139,
  88.785046728971963,
  91.995123933360418,
  77.570093457943926,
  77.529459569280775,
  79.236082893132874,
  85.046728971962608,
  76.391710686712727,
  76.107273466070708,
  80.495733441690362,
  80.373831775700936,
  83.055668427468504,
  81.023973994311262,
  80.86143843965867,
  79.195449004469737,
  76.107273466070708,
  74.197480698902879,
  81.552214546932149,
  75.050792360828936,
  74.806989028850055,
  73.709874034945145,
  74.603819585534339,
  75.294595692807803,
  75.335229581470941,
  80.292563998374646,
  81.227143437626978,
  82.121088988216172,
  78.626574563185699,
  86.021942299878091,
  85.29053230394149,
  79.073547338480282,
  78.464039008533121,
  83.58390898008939
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
